# Update Efnb2-Ephb3 sheet (NatmiData LR-pairs) with recalculated TPM-derived values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 45.71598933333333
    "H2" = 137.147968
    "I2" = 0.6549002937372808
    "J2" = 0.6549002937372808
    "M2" = 0.1709536666666667
    "N2" = 0.512861
    "O2" = 0.007882947722998253
    "P2" = 0.007882947722998253
    "Q2" = 7.815316001827556
    "R2" = 70.337844016448
    "S2" = 0.005162544779307185
    "T2" = 0.005162544779307185
    "G3" = 45.71598933333333
    "H3" = 137.147968
    "I3" = 0.6549002937372808
    "J3" = 0.6549002937372808
    "O3" = 0.7927950496303802
    "P3" = 0.7927950496303802
    "Q3" = 785.9932674004052
    "R3" = 7073.939406603648
    "S3" = 0.5192017108763981
    "T3" = 0.5192017108763981
    "G4" = 45.71598933333333
    "H4" = 137.147968
    "I4" = 0.6549002937372808
    "J4" = 0.6549002937372808
    "M4" = 4.322599666666666
    "N4" = 12.967799
    "O4" = 0.1993220026466216
    "P4" = 0.1993220026466216
    "Q4" = 197.6119202536035
    "R4" = 1778.507282282432
    "S4" = 0.1305360380815755
    "T4" = 0.1305360380815755
    "I5" = 0.1818108415648851
    "J5" = 0.1818108415648851
    "M5" = 0.1709536666666667
    "N5" = 0.512861
    "O5" = 0.007882947722998253
    "P5" = 0.007882947722998253
    "Q5" = 2.169657263824333
    "R5" = 19.526915374419
    "S5" = 0.001433205359530307
    "T5" = 0.001433205359530307
    "I6" = 0.1818108415648851
    "J6" = 0.1818108415648851
    "O6" = 0.7927950496303802
    "P6" = 0.7927950496303802
    "S6" = 0.1441387351617743
    "T6" = 0.1441387351617742
    "I7" = 0.1818108415648851
    "J7" = 0.1818108415648851
    "M7" = 4.322599666666666
    "N7" = 12.967799
    "O7" = 0.1993220026466216
    "P7" = 0.1993220026466216
    "Q7" = 54.86024341130233
    "R7" = 493.7421907017209
    "S7" = 0.03623890104358053
    "T7" = 0.03623890104358053
    "G8" = 11.24784666666667
    "H8" = 33.74354
    "I8" = 0.161130015850732
    "J8" = 0.161130015850732
    "M8" = 0.1709536666666667
    "N8" = 0.512861
    "O8" = 0.007882947722998253
    "P8" = 0.007882947722998253
    "Q8" = 1.922860629771111
    "R8" = 17.30574566794
    "S8" = 0.0012701794915572
    "T8" = 0.0012701794915572
    "G9" = 11.24784666666667
    "H9" = 33.74354
    "I9" = 0.161130015850732
    "J9" = 0.161130015850732
    "O9" = 0.7927950496303802
    "P9" = 0.7927950496303802
    "Q9" = 193.3838003218266
    "R9" = 1740.45420289644
    "S9" = 0.127743078913325
    "T9" = 0.127743078913325
    "G10" = 11.24784666666667
    "H10" = 33.74354
    "I10" = 0.161130015850732
    "J10" = 0.161130015850732
    "M10" = 4.322599666666666
    "N10" = 12.967799
    "O10" = 0.1993220026466216
    "P10" = 0.1993220026466216
    "Q10" = 48.6199382520511
    "R10" = 437.57944426846
    "S10" = 0.03211675744584978
    "T10" = 0.03211675744584978
    "E11" = 2.0
    "F11" = 0.6666666666666666
    "G11" = 0.1507006666666667
    "H11" = 0.452102
    "I11" = 0.00215884884710222
    "J11" = 0.00215884884710222
    "M11" = 0.1709536666666667
    "N11" = 0.512861
    "O11" = 0.007882947722998253
    "P11" = 0.007882947722998253
    "Q11" = 0.02576283153577778
    "R11" = 0.231865483822
    "S11" = 0.00001701809260356185
    "T11" = 0.00001701809260356184
    "E12" = 2.0
    "F12" = 0.6666666666666666
    "G12" = 0.1507006666666667
    "H12" = 0.452102
    "I12" = 0.00215884884710222
    "J12" = 0.00215884884710222
    "O12" = 0.7927950496303802
    "P12" = 0.7927950496303802
    "Q12" = 2.590990835374666
    "R12" = 23.318917518372
    "S12" = 0.001711524678882894
    "T12" = 0.001711524678882893
    "E13" = 2.0
    "F13" = 0.6666666666666666
    "G13" = 0.1507006666666667
    "H13" = 0.452102
    "I13" = 0.00215884884710222
    "J13" = 0.00215884884710222
    "M13" = 4.322599666666666
    "N13" = 12.967799
    "O13" = 0.1993220026466216
    "P13" = 0.1993220026466216
    "Q13" = 0.6514186514997777
    "R13" = 5.862767863498
    "S13" = 0.0004303060756157647
    "T13" = 0.0004303060756157647
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
